$wb = $excel.ActiveWorkbook

# Rename sheet "Fortinet_1" to "Fortinet1"
$sheetFortinet = $wb.Worksheets.Item("Fortinet_1")
$sheetFortinet.Name = "Fortinet1"

$sheetCpmgmt = $wb.Worksheets.Item("CPMGMT")

# Update B2 cell value on Fortinet1 sheet from "Internal DB" to "Internal_DB"
$sheetFortinet.Range("B2").Value = "Internal_DB"

# Update B5 cell value on CPMGMT sheet from "Fortinet_1" to "Fortinet1" (matches renamed sheet)
$sheetCpmgmt.Range("B5").Value = "Fortinet1"

# On CPMGMT sheet, set selection to B5 (not active anymore)
$sheetCpmgmt.Activate()
$sheetCpmgmt.Range("B5").Select()

# Make Fortinet1 the active/selected sheet, with selection B9
$sheetFortinet.Activate()
$sheetFortinet.Range("B9").Select()
